$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add the new "Any (zi-)" column (G) values
$ws.Range("G2").Value = "Any (zi-)"
$ws.Range("G3").Value = "zi-tetl"
$ws.Range("G4").Value = "zi-ad"
$ws.Range("G5").Value = "zi-zhia"
$ws.Range("G6").Value = "zi-qez"
$ws.Range("G7").Value = "zi-patle"
$ws.Range("G8").Value = "zi-qenta"
$ws.Range("G9").Value = "zi-qik"
$ws.Range("G10").Value = "zi-adl"
$ws.Range("G11").Value = "zi-ochti"

# 2. Give column G (rows 2-11) the same formatting as column F (header font/border
#    for row 2, plain bordered font for rows 3-11)
$ws.Range("F2:F11").Copy()
$ws.Range("G2:G11").PasteSpecial(-4122)

# 3. Re-merge the title row across the new column
$ws.Range("A1:F1").UnMerge()
$ws.Range("A1:G1").Merge()

# 4. Extend the title row's bordered look to the new column (copy border-only
#    look from a bordered cell, then reapply the title's font/alignment)
$ws.Range("A2").Copy()
$ws.Range("A1:G1").PasteSpecial(-4122)
$titleRng = $ws.Range("A1:G1")
$titleRng.Font.Bold = $true
$titleRng.Font.Size = 18
$titleRng.HorizontalAlignment = -4108

# 5. Update selection to reflect where the user ended up after editing
$ws.Range("G14").Select()

Write-Host "Edit complete"
